# Auto-generated edit: refresh cryptocurrency price/volume figures.
# (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A few "Price" cells are plain decimal numbers (e.g. "206.29"). The sheet
# stores every Price/Volume cell as literal text, so mark those specific
# cells as Text first -- otherwise Excel would auto-convert the numeric-
# looking string into a real number (and could drop a trailing zero).
$numericLookingCells = @("D5", "D9", "D10", "D17", "D19", "D23", "D24", "D25", "D26", "D32", "D33", "D36", "D39", "D40", "D43", "D44", "D46", "D48", "D51")
foreach ($cellRef in $numericLookingCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = "26.925.38"
$ws.Range("E2").Value = "  -0.39%  "
$ws.Range("D3").Value = "1.551.76"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  -0.34%  "
$ws.Range("D5").Value = "206.29"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("E8").Value = "  +1.27%  "
$ws.Range("D9").Value = "0.248"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").Value = "0.0594"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D12").Value = "1.771.17"
$ws.Range("E12").Value = "  -0.47%  "
$ws.Range("D13").Value = "1.537.26"
$ws.Range("E13").Value = "  -1.26%  "
$ws.Range("E14").Value = "  +0.32%  "
$ws.Range("D16").Value = "26.907.06"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "61.58"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("E18").Value = "  +3.29%  "
$ws.Range("D19").Value = "216.73"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("E22").Value = "  +1.01%  "
$ws.Range("D23").Value = "9.18"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").Value = "1.94"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").Value = "153.50"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").Value = "6.65"
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +0.49%  "
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("E31").Value = "  -1.18%  "
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -0.40%  "
$ws.Range("D33").Value = "3.10"
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("D34").Value = "1.406.51"
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("E35").Value = "  +1.86%  "
$ws.Range("D36").Value = "0.963"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  +0.09%  "
$ws.Range("D39").Value = "0.526"
$ws.Range("E39").Value = "  +0.36%  "
$ws.Range("D40").Value = "0.806"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("E42").Value = "  +3.47%  "
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  +0.75%  "
$ws.Range("D44").Value = "2.29"
$ws.Range("E44").Value = "  +1.08%  "
$ws.Range("E45").Value = "  +0.42%  "
$ws.Range("D46").Value = "1.73"
$ws.Range("E46").Value = "  -1.72%  "
$ws.Range("D47").Value = "1.685.04"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "87.29"
$ws.Range("E48").Value = "  +1.11%  "
$ws.Range("E49").Value = "  +1.24%  "
$ws.Range("D50").Value = "0.0₆0100"
$ws.Range("E50").Value = "  +3.19%  "
$ws.Range("D51").Value = "0.0956"
$ws.Range("E51").Value = "  -0.34%  "
